$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 content updates: blog series shifted from 113/114/115 to 114/115/116,
# and the meetup "Details" link was updated to a new Meetup.com event id.

$ws.Range("B7").Value = "type: featured_blog`nwidth: 2`nheight: 1`nh3: Rules of being a good desi`np: Some rules to follow if you want to lit Pakistan brighter. We here at zakatlists are bounded by these rules. 😀`ndate: 6 Apr 2020`nauthor: <a href=https://justaashir.com target=_blank>Aashir</a>"

$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 116"

$ws.Range("D7").Value = "type: meetup`nwidth: 2`nheight: 1`nh3: Meetup coming in`ndate: 2020,4,3,10,30,0,0`nbutton.default: Speak*goto(`"https://forms.gle/dyydXFRSsKzeH4hZ6`")`nbutton.default: Attend*goto(`"https://youtu.be/vscn-HP932E`")`nbutton.default: Details*goto(`"https://www.meetup.com/techshek/events/270179438/`")"

$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 115"

$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 114"

# Selection moved from I7 to D7
$ws.Range("D7").Select() | Out-Null
